$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "continuous frames" annotation block in columns G:J (rows 1-11),
#    written in row-major order so shared strings are appended in the same
#    sequence as the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "0-570"
$ws.Range("H1").Value = "670-13180"
$ws.Range("I1").Value = "13270-29140"
$ws.Range("J1").Value = "29210-end"

$ws.Range("G2").Value = "0-760"
$ws.Range("H2").Value = "880-16370"
$ws.Range("I2").Value = "16500-29200"
$ws.Range("J2").Value = "29330-end"

$ws.Range("G3").Value = "0-11370"
$ws.Range("H3").Value = "11420-22370"
$ws.Range("I3").Value = "22590-end"

$ws.Range("G4").Value = "0-15900"
$ws.Range("H4").Value = "16030-23940"
$ws.Range("I4").Value = "24090-end"

$ws.Range("G5").Value = "0-2120"
$ws.Range("H5").Value = "2180-25010"
$ws.Range("I5").Value = "25040-end"

$ws.Range("G6").Value = "0-750"
$ws.Range("H6").Value = "810-12580"
$ws.Range("I6").Value = "12660-23540"
$ws.Range("J6").Value = "23650-end"

$ws.Range("G7").Value = "0-2510"
$ws.Range("H7").Value = "2550-20870"
$ws.Range("I7").Value = "20950-end"

$ws.Range("G8").Value = "0-1250"
$ws.Range("H8").Value = "1300-12570"
$ws.Range("I8").Value = "12680-24520"
$ws.Range("J8").Value = "24640-end"

$ws.Range("G9").Value = "0-4310"
$ws.Range("H9").Value = "4370-end"

$ws.Range("G10").Value = "0-1170"
$ws.Range("H10").Value = "1250-25570"
$ws.Range("I10").Value = "25700-28260"
$ws.Range("J10").Value = "28360-end"

$ws.Range("G11").Value = "0-20680"
$ws.Range("H11").Value = "20740-28400"
$ws.Range("I11").Value = "28530-end"

# Extra helper numbers next to row 1
$ws.Range("L1").Value = 29140
$ws.Range("M1").Value = 29210
$ws.Range("N1").Value = 32399

# ---------------------------------------------------------------------------
# 2) Formatting: two small Arial-10 font variants are introduced.
#    Style "A" (theme text colour) covers the whole G:J annotation block
#    (including the blank trailing cells), style "B" (explicit black) is
#    then re-applied on top of it for the populated cells of rows 9-11.
#    A single seed cell is formatted directly and its format is fanned out
#    with Copy/PasteSpecial so no extra transient font/style entries leak
#    into styles.xml.
# ---------------------------------------------------------------------------
$seedA = $ws.Range("G1")
$seedA.Font.Name = "Arial"
$seedA.Font.Size = 10

$styleARange = $ws.Range("G1:J8,I9:J9,J11")
foreach ($area in $styleARange.Areas) {
    $seedA.Copy()
    $area.PasteSpecial(-4122)
}

$seedB = $ws.Range("G9")
$seedB.Font.Name = "Arial"
$seedB.Font.Size = 10
$seedB.Font.Color = 0

$styleBRange = $ws.Range("G9:H9,G10:J10,G11:I11")
foreach ($area in $styleBRange.Areas) {
    $seedB.Copy()
    $area.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) View state: move the active selection from A48 to A12 (and drop the
#    stale "topLeftCell" scroll position in the process).
# ---------------------------------------------------------------------------
$ws.Range("A12").Select()

# ---------------------------------------------------------------------------
# 4) Window placement metadata.
# ---------------------------------------------------------------------------
try { $excel.Left = 1000 } catch {}
try { $excel.Top = 60 } catch {}
try { $excel.ActiveWindow.Left = 1000 } catch {}
try { $excel.ActiveWindow.Top = 60 } catch {}
